$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chapter 1 finished -> row 17's Number goes from 15 to 16
$ws.Cells.Item(17, 1).Value = 16

# Remove the stray explicit "applied number format" styling on the Hours
# column (D5:D17) that no longer serves a purpose.
$ws.Range("D5:D17").ClearFormats()

# Chapter 2 started -> append a new entry row (row 18)
$ws.Cells.Item(18, 1).Value = 17

# Column B holds dates stored as plain text (e.g. "2019.02.19"), so force
# text formatting before assigning the value to avoid Excel's automatic
# date parsing, then clear the formatting again so no style is retained.
$ws.Cells.Item(18, 2).NumberFormat = "@"
$ws.Cells.Item(18, 2).Value = "2019.03.17"
$ws.Cells.Item(18, 2).ClearFormats()

$ws.Cells.Item(18, 3).Value = "#latex"
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 5).Value = "Hub calculations"

# Reflect the last place the user clicked after entering the new row.
$ws.Range("E19").Select()
